$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend header sequence with P1=14, Q1=15 (same style as O1) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$hdr = $ws.Range("P1:Q1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# --- Rows 2-25: swap I/K and M/O values, and add new P,Q columns = 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column = 2
}

Write-Host "done"
